$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C22").Value = "'322"
$ws.Range("D22").Value = "'935504.07"
$ws.Range("C22:D22").Style = "Normal"

$ws.Range("C23").Value = "'115"
$ws.Range("D23").Value = "'458348.00"
$ws.Range("C23:D23").Style = "Normal"

$ws.Range("C25").Value = "'7"
$ws.Range("D25").Value = "'57000.00"
$ws.Range("C25:D25").Style = "Normal"

$ws.Range("C33").Value = "'93"
$ws.Range("D33").Value = "'243826.00"
$ws.Range("C33:D33").Style = "Normal"

$ws.Range("C34").Value = "'508"
$ws.Range("D34").Value = "'1563122.82"
$ws.Range("C34:D34").Style = "Normal"

$ws.Range("C35").Value = "'205"
$ws.Range("D35").Value = "'990347.11"
$ws.Range("C35:D35").Style = "Normal"

$ws.Range("C38").Value = "'19"
$ws.Range("D38").Value = "'42200.00"
$ws.Range("C38:D38").Style = "Normal"

$ws.Range("C45").Value = "'25"
$ws.Range("D45").Value = "'98621.84"
$ws.Range("C45:D45").Style = "Normal"

$ws.Range("C46").Value = "'68"
$ws.Range("D46").Value = "'287139.82"
$ws.Range("C46:D46").Style = "Normal"

$ws.Range("C47").Value = "'40"
$ws.Range("D47").Value = "'231937.00"
$ws.Range("C47:D47").Style = "Normal"

$ws.Range("C49").Value = "'6"
$ws.Range("D49").Value = "'15850.00"
$ws.Range("C49:D49").Style = "Normal"

$ws.Range("C50").Value = "'93"
$ws.Range("D50").Value = "'257768.17"
$ws.Range("C50:D50").Style = "Normal"

$ws.Range("C51").Value = "'541"
$ws.Range("D51").Value = "'1774758.52"
$ws.Range("C51:D51").Style = "Normal"

$ws.Range("C52").Value = "'250"
$ws.Range("D52").Value = "'1045235.76"
$ws.Range("C52:D52").Style = "Normal"

$ws.Range("C53").Value = "'82"
$ws.Range("D53").Value = "'466378.23"
$ws.Range("C53:D53").Style = "Normal"

$ws.Range("C55").Value = "'16"
$ws.Range("D55").Value = "'44720.65"
$ws.Range("C55:D55").Style = "Normal"

$ws.Range("C56").Value = "'668"
$ws.Range("D56").Value = "'1670796.41"
$ws.Range("C56:D56").Style = "Normal"

$ws.Range("C57").Value = "'3285"
$ws.Range("D57").Value = "'9856825.10"
$ws.Range("C57:D57").Style = "Normal"

$ws.Range("C58").Value = "'1695"
$ws.Range("D58").Value = "'6745611.92"
$ws.Range("C58:D58").Style = "Normal"

$ws.Range("C59").Value = "'578"
$ws.Range("D59").Value = "'2717640.96"
$ws.Range("C59:D59").Style = "Normal"

$ws.Range("C79").Value = "'220"
$ws.Range("D79").Value = "'557826.09"
$ws.Range("C79:D79").Style = "Normal"

$ws.Range("C80").Value = "'849"
$ws.Range("D80").Value = "'2613791.11"
$ws.Range("C80:D80").Style = "Normal"

$ws.Range("C84").Value = "'30"
$ws.Range("D84").Value = "'61500.00"
$ws.Range("C84:D84").Style = "Normal"

$ws.Range("C97").Value = "'269"
$ws.Range("D97").Value = "'695896.75"
$ws.Range("C97:D97").Style = "Normal"

$ws.Range("C98").Value = "'1167"
$ws.Range("D98").Value = "'3478927.44"
$ws.Range("C98:D98").Style = "Normal"

$ws.Range("C99").Value = "'437"
$ws.Range("D99").Value = "'1712594.02"
$ws.Range("C99:D99").Style = "Normal"

$ws.Range("C100").Value = "'115"
$ws.Range("D100").Value = "'517000.00"
$ws.Range("C100:D100").Style = "Normal"

$ws.Range("C101").Value = "'30"
$ws.Range("D101").Value = "'179157.00"
$ws.Range("C101:D101").Style = "Normal"

$ws.Range("C102").Value = "'63"
$ws.Range("D102").Value = "'137000.00"
$ws.Range("C102:D102").Style = "Normal"
